$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-13 Friday" "2024-09-14 Saturday"

Replace-Text "534×9=" "378×9="
Replace-Text "337×9=" "633×5="
Replace-Text "871×6=" "787×9="
Replace-Text "231×3=" "207×8="
Replace-Text "249×3=" "805×3="
Replace-Text "112×7=" "495×9="
Replace-Text "170×5=" "922×8="
Replace-Text "643×4=" "136×8="
Replace-Text "276×4=" "982×9="
Replace-Text "925×4=" "206×8="
Replace-Text "706×4=" "459×7="
Replace-Text "594×3=" "271×8="
Replace-Text "395×7=" "361×3="
Replace-Text "740×5=" "752×2="
Replace-Text "376×8=" "457×4="
Replace-Text "653×8=" "871×9="
Replace-Text "255×2=" "310×8="
Replace-Text "977×9=" "762×2="
Replace-Text "174×3=" "117×9="
Replace-Text "384×8=" "455×9="
Replace-Text "832×4=" "445×5="
Replace-Text "216×8=" "706×8="
Replace-Text "196×3=" "905×9="
Replace-Text "107×3=" "491×4="
Replace-Text "256×4=" "356×4="
